$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to retain text formatting so purely numeric-looking
# values (e.g. "1.005") are not converted into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.420.46'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.791.97'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +1.03%  '
$ws.Range("D5").Value = '338.19'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = '0.3792'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").Value = '0.3452'
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("D9").Value = '48.75'
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("D10").Value = '1.199'
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").Value = '0.07511'
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '21.92'
$ws.Range("E13").Value = '  +6.07%  '
$ws.Range("D14").Value = '6.468'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '1.791.43'
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").Value = '0.00001102'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '0.06664'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").Value = '84.73'
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '6.530'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  +2.56%  '
$ws.Range("D23").Value = '27.412.81'
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").Value = '12.52'
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("D25").Value = '2.439'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '2.566'
$ws.Range("E26").Value = '  +4.50%  '
$ws.Range("D27").Value = '1.499'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '21.46'
$ws.Range("E28").Value = '  +8.38%  '
$ws.Range("D29").Value = '153.43'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '1.994.27'
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").Value = '133.30'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").Value = '6.092'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").Value = '0.08685'
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").Value = '13.20'
$ws.Range("E35").Value = '  +1.11%  '
$ws.Range("D36").Value = '1.657'
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = '0.6891'
$ws.Range("E38").Value = '  +8.09%  '
$ws.Range("D39").Value = '0.06375'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").Value = '8.858'
$ws.Range("E40").Value = '  +2.87%  '
$ws.Range("D41").Value = '0.2200'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '0.02347'
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").Value = '1.263'
$ws.Range("E43").Value = '  +2.82%  '
$ws.Range("D44").Value = '14.36'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").Value = '0.6426'
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("D47").Value = '3.869'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").Value = '2.135'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").Value = '129.52'
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("D50").Value = '0.07196'
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").Value = '79.34'
$ws.Range("E51").Value = '  +0.54%  '

# Restore the default cell style (the diff does not alter styling),
# while keeping the values stored as text.
$ws.Range("D2:D51").Style = "Normal"

